# Update Daily Report: 2026-02-26
# - Append a new day's depository data (date serial 46078) to Daily_Data
# - Refresh Today_Summary (latest-day snapshot) with the new day's figures
# - Refresh Monthly_Stats (month-to-date rollups) to include the new day

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Daily_Data: append rows 290-313 for date serial 46078
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(46078, 'ASAHI DEPOSITORY LLC Registered', 23291615.992, 0, 0, 0, 0, 23291615.992),
    @(46078, 'ASAHI DEPOSITORY LLC Eligible', 1703429.248, 0, 0, 0, 0, 1703429.248),
    @(46078, 'BRINK''S, INC. Registered', 14694890.273, 0, 0, 0, -217756.88, 14477133.393),
    @(46078, 'BRINK''S, INC. Eligible', 40424764.88, 0, 0, 0, 217756.88, 40642521.76),
    @(46078, 'CNT DEPOSITORY, INC. Registered', 12170205.469, 0, 0, 0, 0, 12170205.469),
    @(46078, 'CNT DEPOSITORY, INC. Eligible', 13861333.923, 0, 0, 0, 0, 13861333.923),
    @(46078, 'DELAWARE DEPOSITORY Registered', 1532776.423, 0, 0, 0, 0, 1532776.423),
    @(46078, 'DELAWARE DEPOSITORY Eligible', 15761931.149, 0, 3158.747, -3158.747, 0, 15758772.402),
    @(46078, 'HSBC BANK, USA Registered', 3387219.03, 0, 0, 0, 0, 3387219.03),
    @(46078, 'HSBC BANK, USA Eligible', 18388203.313, 0, 623546.3, -623546.3, 0, 17764657.013),
    @(46078, 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 273789.87, 0, 0, 0, 0, 273789.87),
    @(46078, 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 3295246.644, 0, 0, 0, 0, 3295246.644),
    @(46078, 'JP MORGAN CHASE BANK NA Registered', 12000343.77, 0, 0, 0, 0, 12000343.77),
    @(46078, 'JP MORGAN CHASE BANK NA Eligible', 142516268.453, 0, 24904.77, -24904.77, 0, 142491363.683),
    @(46078, 'LOOMIS INTERNATIONAL (US) LLC Registered', 6842629.447, 0, 0, 0, 0, 6842629.447),
    @(46078, 'LOOMIS INTERNATIONAL (US) LLC Eligible', 23512931.636, 0, 0, 0, 0, 23512931.636),
    @(46078, 'MALCA-AMIT ARMORED, INC. Registered', 0, 0, 0, 0, 0, 0),
    @(46078, 'MALCA-AMIT ARMORED, INC. Eligible', 0, 0, 0, 0, 0, 0),
    @(46078, 'MALCA-AMIT USA, LLC Registered', 949634.064, 0, 0, 0, 0, 949634.064),
    @(46078, 'MALCA-AMIT USA, LLC Eligible', 1073898.377, 0, 0, 0, 0, 1073898.377),
    @(46078, 'MANFRA, TORDELLA & BROOKES, LLC Registered', 5871594.333, 0, 0, 0, 0, 5871594.333),
    @(46078, 'MANFRA, TORDELLA & BROOKES, LLC Eligible', 12604051.607, 0, 619157.709, -619157.709, 0, 11984893.898),
    @(46078, 'STONEX PRECIOUS METALS LLC Registered', 5268155.68, 65301.69, 0, 65301.69, 0, 5333457.37),
    @(46078, 'STONEX PRECIOUS METALS LLC Eligible', 2419487.69, 0, 0, 0, 0, 2419487.69)
)

$r = 290
foreach ($row in $newRows) {
    $wsData.Cells.Item($r, 1).Value = $row[0]
    $wsData.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsData.Cells.Item($r, 2).Value = $row[1]
    $wsData.Cells.Item($r, 3).Value = $row[2]
    $wsData.Cells.Item($r, 4).Value = $row[3]
    $wsData.Cells.Item($r, 5).Value = $row[4]
    $wsData.Cells.Item($r, 6).Value = $row[5]
    $wsData.Cells.Item($r, 7).Value = $row[6]
    $wsData.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Today_Summary: depositories whose Eligible/Registered split changed
#    on the new day (totals per depository)
# ---------------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("Today_Summary")

# BRINK'S, INC. (row 3)
$wsToday.Cells.Item(3, 2).Value = 40642521.76
$wsToday.Cells.Item(3, 3).Value = 14477133.393

# DELAWARE DEPOSITORY (row 5)
$wsToday.Cells.Item(5, 2).Value = 15758772.402
$wsToday.Cells.Item(5, 4).Value = 17291548.825

# HSBC BANK, USA (row 6)
$wsToday.Cells.Item(6, 2).Value = 17764657.013
$wsToday.Cells.Item(6, 4).Value = 21151876.043

# JP MORGAN CHASE BANK NA (row 8)
$wsToday.Cells.Item(8, 2).Value = 142491363.683
$wsToday.Cells.Item(8, 4).Value = 154491707.453

# MANFRA, TORDELLA & BROOKES, LLC (row 12)
$wsToday.Cells.Item(12, 2).Value = 11984893.898
$wsToday.Cells.Item(12, 4).Value = 17856488.231

# STONEX PRECIOUS METALS LLC (row 13)
$wsToday.Cells.Item(13, 3).Value = 5333457.37
$wsToday.Cells.Item(13, 4).Value = 7752945.060000001

# ---------------------------------------------------------------------------
# 3) Monthly_Stats: month-to-date rollups now include the new day
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# Grand total summary row (row 2)
$wsMonthly.Cells.Item(2, 2).Value = 274508536.274
$wsMonthly.Cells.Item(2, 3).Value = 86130399.161
$wsMonthly.Cells.Item(2, 4).Value = 360638935.435

# Detail rows (YearMonth / Region_Type / RECEIVED / WITHDRAWN / TOTAL_TODAY)
# BRINK'S, INC. Eligible (row 9)
$wsMonthly.Cells.Item(9, 5).Value = 40642521.76

# BRINK'S, INC. Registered (row 10)
$wsMonthly.Cells.Item(10, 5).Value = 14477133.393

# DELAWARE DEPOSITORY Eligible (row 13)
$wsMonthly.Cells.Item(13, 4).Value = 704523.667
$wsMonthly.Cells.Item(13, 5).Value = 15758772.402

# HSBC BANK, USA Eligible (row 15)
$wsMonthly.Cells.Item(15, 4).Value = 3581337.69
$wsMonthly.Cells.Item(15, 5).Value = 17764657.013

# JP MORGAN CHASE BANK NA Eligible (row 19)
$wsMonthly.Cells.Item(19, 4).Value = 14372437.75
$wsMonthly.Cells.Item(19, 5).Value = 142491363.683

# MANFRA, TORDELLA & BROOKES, LLC Eligible (row 27)
$wsMonthly.Cells.Item(27, 4).Value = 1661718.932
$wsMonthly.Cells.Item(27, 5).Value = 11984893.898

# STONEX PRECIOUS METALS LLC Registered (row 30)
$wsMonthly.Cells.Item(30, 3).Value = 85064.87
$wsMonthly.Cells.Item(30, 5).Value = 5333457.37
